# Weekly price-sheet update: a new weekly observation (2021-11-09, serial
# 44509) is inserted right after the current row 161, pushing the existing
# rows 161-167 down to 162-168 (dimension grows from R167 to R168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 161; existing data (rows 161-167) shifts to 162-168.
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new weekly observation.
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44509
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100112003
$ws.Range("G161").Value = "Ajo"
$ws.Range("H161").Value = "Chino"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 480
$ws.Range("K161").Value = 19000
$ws.Range("L161").Value = 19500
$ws.Range("M161").Value = 19250
$ws.Range("N161").Value = "`$/caja 10 kilos"
$ws.Range("O161").Value = "China"
$ws.Range("P161").Value = 1925
$ws.Range("Q161").Value = 10
$ws.Range("R161").Value = "Hortaliza"
